$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 108, shifting existing rows 108:167 down to 109:168
$ws.Rows.Item(108).Insert()

# Populate the newly inserted row 108 with the new data record
$ws.Range("A108").Value = 10
$ws.Range("B108").Value = "Vega Modelo de Temuco"
$ws.Range("C108").Value = "La Araucanía"
$ws.Range("D108").Value = 44572
$ws.Range("E108").Value = 9
$ws.Range("F108").Value = 100112043
$ws.Range("G108").Value = "Pepino dulce"
$ws.Range("H108").Value = "Cultivar IV Región"
$ws.Range("I108").Value = "Primera"
$ws.Range("J108").Value = 210
$ws.Range("K108").Value = 25000
$ws.Range("L108").Value = 25000
$ws.Range("M108").Value = 25000
$ws.Range("N108").Value = "$/bandeja 18 kilos"
$ws.Range("O108").Value = "Provincia de Limarí"
$ws.Range("P108").Value = 1389
$ws.Range("Q108").Value = 18
$ws.Range("R108").Value = "Hortaliza"
